# FN-3460: fix numerical rounding errors during report upload
# Adds two new facility-utilisation rows ("Crumpet" and "Scone") to the
# next-week utilisation report fixture, following the same layout/style
# as the existing rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 carries the formatting we want row 5 to inherit (row 5 was a
# completely blank row in the template, so copy the style down first).
$ws.Range("A4:K4").Copy()
$ws.Range("A5:K5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 5: Crumpet GEF / Crumpet exporter ---
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"
$ws.Range("K5").Clear()

# --- Row 6: Scone GEF / Scone exporter ---
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Match the author's final selection in the saved file.
$ws.Range("A5:J6").Select()
